$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H93").Value = 45375.5
$ws.Range("J93").Value = 45375.5
$ws.Range("L93").Value = 45375.5
$ws.Range("N93").Value = -50367.5

$ws.Range("H98").Value = 3522.795
$ws.Range("J98").Value = 6211.2
$ws.Range("L98").Value = 6211.2
$ws.Range("N98").Value = -9207.200000000001

$ws.Range("H122").Value = 3522.795
$ws.Range("J122").Value = 6211.2
$ws.Range("L122").Value = 18633.6
$ws.Range("N122").Value = -23533.6

$ws.Range("H129").Value = 1023.60345
$ws.Range("J129").Value = 1257.2325
$ws.Range("L129").Value = 3771.6975
$ws.Range("N129").Value = -13771.6975


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 11370.418
$ws.Range("I32").Value = 11158.3
$ws.Range("J32").Value = 13491.6
$ws.Range("K32").Value = 11158.3
$ws.Range("L32").Value = 13491.6
$ws.Range("M32").Value = -10871.3
$ws.Range("N32").Value = -14065.6

$ws.Range("H61").Value = 2437.48
$ws.Range("I61").Value = 2329.7334
$ws.Range("J61").Value = 2599.1
$ws.Range("K61").Value = 2329.7334
$ws.Range("L61").Value = 2599.1
$ws.Range("M61").Value = -2117.7334
$ws.Range("N61").Value = -3023.1

$ws.Range("H74").Value = 929.8929000000001
$ws.Range("I74").Value = 870.5789
$ws.Range("J74").Value = 1055.1111
$ws.Range("K74").Value = 870.5789
$ws.Range("L74").Value = 1055.1111
$ws.Range("M74").Value = 3.421100000000024
$ws.Range("N74").Value = -2803.1111

$ws.Range("H77").Value = 929.8929000000001
$ws.Range("I77").Value = 870.5789
$ws.Range("J77").Value = 1055.1111
$ws.Range("K77").Value = 4352.8945
$ws.Range("L77").Value = 5275.5555
$ws.Range("M77").Value = 15.10549999999967
$ws.Range("N77").Value = -14011.5555

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 350000
$ws.Range("J92").Value = 350000
$ws.Range("L92").Value = 350000
$ws.Range("N92").Value = -354992

$ws.Range("H128").Value = 45000
$ws.Range("J128").Value = 45000
$ws.Range("L128").Value = 45000
$ws.Range("N128").Value = -54960

$ws.Range("H132").Value = 4553.268
$ws.Range("I132").Value = 5340.393
$ws.Range("J132").Value = 2857.923
$ws.Range("K132").Value = 16021.179
$ws.Range("L132").Value = 8573.769
$ws.Range("M132").Value = -13491.179
$ws.Range("N132").Value = -13633.769

$ws.Range("H136").Value = 2437.48
$ws.Range("I136").Value = 2329.7334
$ws.Range("J136").Value = 2599.1
$ws.Range("K136").Value = 6989.2002
$ws.Range("L136").Value = 7797.299999999999
$ws.Range("M136").Value = -4439.2002
$ws.Range("N136").Value = -12897.3


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H20").Value = 36584.793
$ws.Range("I20").Value = 53926.42
$ws.Range("J20").Value = 3635.7
$ws.Range("K20").Value = 53926.42
$ws.Range("L20").Value = 3635.7
$ws.Range("M20").Value = -53679.42
$ws.Range("N20").Value = -4129.7

$ws.Range("H94").Value = 336333.34
$ws.Range("J94").Value = 336333.34
$ws.Range("L94").Value = 336333.34
$ws.Range("N94").Value = -337235.34

$ws.Range("H99").Value = 1397.9524
$ws.Range("I99").Value = 1135.5
$ws.Range("J99").Value = 1922.8572
$ws.Range("K99").Value = 1135.5
$ws.Range("L99").Value = 1922.8572
$ws.Range("M99").Value = 362.5
$ws.Range("N99").Value = -4918.8572

$ws.Range("H105").Value = 3569.8572
$ws.Range("I105").Value = 3333.1667
$ws.Range("K105").Value = 3333.1667
$ws.Range("M105").Value = -1586.1667

$ws.Range("H134").Value = 2560.625
$ws.Range("I134").Value = 2139.2917
$ws.Range("J134").Value = 3824.625
$ws.Range("K134").Value = 6417.875100000001
$ws.Range("L134").Value = 11473.875
$ws.Range("M134").Value = -3882.875100000001
$ws.Range("N134").Value = -16543.875


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2041.4286
$ws.Range("I31").Value = 1506.1538
$ws.Range("J31").Value = 9000
$ws.Range("K31").Value = 1506.1538
$ws.Range("L31").Value = 9000
$ws.Range("M31").Value = -1211.1538
$ws.Range("N31").Value = -9590

$ws.Range("H34").Value = 2041.4286
$ws.Range("I34").Value = 1506.1538
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 1506.1538
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -1304.1538
$ws.Range("N34").Value = -9404

$ws.Range("H58").Value = 2472153.8
$ws.Range("I58").Value = 3089150.5
$ws.Range("J58").Value = 4166.6665
$ws.Range("K58").Value = 3089150.5
$ws.Range("L58").Value = 4166.6665
$ws.Range("M58").Value = -3088947.5
$ws.Range("N58").Value = -4572.6665

$ws.Range("H132").Value = 331157.16
$ws.Range("I132").Value = 483718.2
$ws.Range("J132").Value = 2564.1538
$ws.Range("K132").Value = 1451154.6
$ws.Range("L132").Value = 7692.4614
$ws.Range("M132").Value = -1448624.6
$ws.Range("N132").Value = -12752.4614

$ws.Range("H134").Value = 1991.6154
$ws.Range("I134").Value = 1408.8422
$ws.Range("J134").Value = 3573.4285
$ws.Range("K134").Value = 4226.5266
$ws.Range("L134").Value = 10720.2855
$ws.Range("M134").Value = -1691.5266
$ws.Range("N134").Value = -15790.2855

$ws.Range("H136").Value = 2472153.8
$ws.Range("I136").Value = 3089150.5
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 9267451.5
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -9264901.5
$ws.Range("N136").Value = -17599.9995


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H32").Value = 1822
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 1822
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 5466
$ws.Range("N32").Value = -6032
$ws.Range("M32").ClearContents()

$ws.Range("H107").Value = 380.41177
$ws.Range("J107").Value = 308.6
$ws.Range("L107").Value = 925.8000000000001
$ws.Range("N107").Value = -4765.8

$ws.Range("H131").Value = 13700139
$ws.Range("I131").Value = 360
$ws.Range("J131").Value = 15153146
$ws.Range("K131").Value = 1080
$ws.Range("L131").Value = 45459438
$ws.Range("M131").Value = 3960
$ws.Range("N131").Value = -45469518

$ws.Range("H132").Value = 1995.6666
$ws.Range("I132").Value = 1593.8667
$ws.Range("K132").Value = 14344.8003
$ws.Range("M132").Value = -11814.8003


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H62").Value = 31000
$ws.Range("J62").Value = 31000
$ws.Range("L62").Value = 31000
$ws.Range("N62").Value = -32372

$ws.Range("H65").Value = 31000
$ws.Range("J65").Value = 31000
$ws.Range("L65").Value = 93000
$ws.Range("N65").Value = -99864

$ws.Range("H70").Value = 5860.0386
$ws.Range("I70").Value = 5663.263
$ws.Range("J70").Value = 6394.143
$ws.Range("K70").Value = 5663.263
$ws.Range("L70").Value = 6394.143
$ws.Range("M70").Value = -5393.263
$ws.Range("N70").Value = -6934.143

$ws.Range("H73").Value = 5860.0386
$ws.Range("I73").Value = 5663.263
$ws.Range("J73").Value = 6394.143
$ws.Range("K73").Value = 5663.263
$ws.Range("L73").Value = 6394.143
$ws.Range("M73").Value = -4727.263
$ws.Range("N73").Value = -8266.143

$ws.Range("H131").Value = 28491.75
$ws.Range("J131").Value = 28491.75
$ws.Range("L131").Value = 28491.75
$ws.Range("N131").Value = -38571.75


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H9").Value = 4250
$ws.Range("I9").Value = 321.7143
$ws.Range("J9").Value = 11124.5
$ws.Range("K9").Value = 321.7143
$ws.Range("L9").Value = 11124.5
$ws.Range("M9").Value = -97.71429999999998
$ws.Range("N9").Value = -11572.5

$ws.Range("H69").Value = 193387.67
$ws.Range("J69").Value = 275081.5
$ws.Range("L69").Value = 275081.5
$ws.Range("N69").Value = -276703.5

$ws.Range("H70").Value = 181721
$ws.Range("J70").Value = 181721
$ws.Range("L70").Value = 181721
$ws.Range("N70").Value = -182261

$ws.Range("H72").Value = 193387.67
$ws.Range("J72").Value = 275081.5
$ws.Range("L72").Value = 825244.5
$ws.Range("N72").Value = -833356.5

$ws.Range("H73").Value = 181721
$ws.Range("J73").Value = 181721
$ws.Range("L73").Value = 181721
$ws.Range("N73").Value = -183593

$ws.Range("H127").Value = 44048.332
$ws.Range("J127").Value = 44048.332
$ws.Range("L127").Value = 44048.332
$ws.Range("N127").Value = -53968.332


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 4333
$ws.Range("I62").Value = 3999.5
$ws.Range("K62").Value = 3999.5
$ws.Range("M62").Value = -3375.5

$ws.Range("H65").Value = 4333
$ws.Range("I65").Value = 3999.5
$ws.Range("K65").Value = 19997.5
$ws.Range("M65").Value = -16877.5

$ws.Range("H126").Value = 3453.6897
$ws.Range("I126").Value = 5169.5557
$ws.Range("J126").Value = 645.9091
$ws.Range("K126").Value = 15508.6671
$ws.Range("L126").Value = 1937.7273
$ws.Range("M126").Value = -13038.6671
$ws.Range("N126").Value = -6877.7273

$ws.Range("H131").Value = 44999.5
$ws.Range("J131").Value = 44999.5
$ws.Range("L131").Value = 44999.5
$ws.Range("N131").Value = -55079.5

